$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of gyroscope data to be inserted at the top of the dataset
# (pushing the existing rows down), and the final row of the old data
# is dropped so the sheet ends at row 21 instead of row 22.
$newRows = @(
    @(-0.2000583708286285, -0.1212567538022995, -0.0207694191485643),
    @(-0.1815796941518783, -0.0572686158120632, 0.08643743395805351),
    @(-0.0739146918058395, -0.1140790879726409, 0.1067487001419067),
    @(-0.0395535230636596, -0.0899499058723449, -0.0404698215425014),
    @(-0.0148134818300604, 0.1036943718791008, -0.1157589629292488),
    @(0.5971207618713379, 1.289536476135254, -0.3637702465057373)
)

$insertCount = $newRows.Count

# Insert blank rows above row 2 to make room for the new data while
# shifting the existing rows down.
$insertRange = $ws.Range("A2:C$(1 + $insertCount)")
$insertRange.Insert()

# The inserted rows pick up the header row's formatting; clear it so the
# new data rows match the plain (unstyled) look of the rest of the data.
$insertRange.ClearFormats()

# Write the new rows into the freshly inserted space.
for ($i = 0; $i -lt $insertCount; $i++) {
    $rowNum = 2 + $i
    $values = $newRows[$i]
    $ws.Cells.Item($rowNum, 1).Value = $values[0]
    $ws.Cells.Item($rowNum, 2).Value = $values[1]
    $ws.Cells.Item($rowNum, 3).Value = $values[2]
}

# The bottom of the original dataset (old rows 16-22, now shifted down
# by $insertCount rows) is no longer part of the dataset; remove those
# rows so the sheet ends at row 21.
$deleteFirst = 16 + $insertCount
$deleteLast = 22 + $insertCount
$ws.Range("A$($deleteFirst):C$($deleteLast)").EntireRow.Delete()
